# Update "想去人数" (number of people interested) counts on the 展览 (Exhibition)
# sheet and the matching rows on the 全部类型 (All types) aggregate sheet.
# 展览!F3: 254 -> 256
# 展览!F4: 894 -> 897
# 展览!F6: 40  -> 41
# 全部类型!F4: 254 -> 256
# 全部类型!F5: 894 -> 897
# 全部类型!F7: 40  -> 41

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 256
$wsExhibit.Range("F4").Value = 897
$wsExhibit.Range("F6").Value = 41

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 256
$wsAll.Range("F5").Value = 897
$wsAll.Range("F7").Value = 41
